# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps written by the handback report job.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
$wsOverview.Range("G2").Value = "2016-08-19 23:09:54"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn.Range("H2").Value = "2016-08-19 23:09:50"
$wsZhCn.Range("K2").Value = "2016-08-19 23:10:15"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsDeDe.Range("H2").Value = "2016-08-19 23:09:54"
$wsDeDe.Range("K2").Value = "2016-08-19 23:10:21"
